$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.214.41'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '1.588.37'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  +1.10%  '
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("D8").Value = "'23.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.81%  '
$ws.Range("D9").Value = "'0.251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = "'0.0599"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("D11").Value = "'0.0889"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.31%  '
$ws.Range("D12").Value = '1.816.09'
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").Value = '1.588.56'
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("E14").Value = '  +0.83%  '
$ws.Range("E15").Value = '  -1.03%  '
$ws.Range("D16").Value = '28.251.38'
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("D17").Value = "'63.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.76%  '
$ws.Range("D18").Value = "'227.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = "'7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  +0.97%  '
$ws.Range("E22").Value = '  -1.95%  '
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").Value = "'151.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.64%  '
$ws.Range("D26").Value = "'15.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.58%  '
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("E28").Value = '  -1.46%  '
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").Value = "'0.0472"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").Value = '1.397.24'
$ws.Range("E34").Value = '  -4.23%  '
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("E36").Value = '  -7.89%  '
$ws.Range("D37").Value = "'2.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.45%  '
$ws.Range("E38").Value = '  -0.77%  '
$ws.Range("E39").Value = '  +5.84%  '
$ws.Range("D40").Value = "'0.540"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").Value = "'0.812"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.93%  '
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("E43").Value = '  -3.52%  '
$ws.Range("E44").Value = '  +3.57%  '
$ws.Range("D45").Value = "'0.979"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("D46").Value = "'64.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("D47").Value = '1.725.62'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("E49").Value = '  +1.87%  '
$ws.Range("E50").Value = '  +7.68%  '
$ws.Range("E51").Value = '  -0.55%  '
